# Weekly update: insert a new price-report row for the current week at
# row 215 of the "Fruta, Vega Monumental Concepción - Piña" sheet. All
# subsequent rows (old 215..293) shift down by one (to 216..294).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the existing row 215 - this pushes rows
# 215..293 down to 216..294 and extends the used range to A1:T294.
$ws.Rows("215:215").Insert()

# Populate the newly inserted row 215 with this week's record.
$ws.Range("A215").Value = 11
$ws.Range("B215").Value = "Vega Monumental Concepción"
$ws.Range("C215").Value = "Bíobío"
$ws.Range("D215").Value = 45146
$ws.Range("E215").Value = 8
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100108
$ws.Range("H215").Value = "Tropicales y subtropicales"
$ws.Range("I215").Value = 100108005
$ws.Range("J215").Value = "Piña"
$ws.Range("K215").Value = "Caramelo"
$ws.Range("L215").Value = "Segunda"
$ws.Range("M215").Value = 250
$ws.Range("N215").Value = 19000
$ws.Range("O215").Value = 20000
$ws.Range("P215").Value = 19400
$ws.Range("Q215").Value = '$/caja 14 unidades'
$ws.Range("R215").Value = "Ecuador"
$ws.Range("S215").Value = 1386
$ws.Range("T215").Value = 14
